# Apply the "Saldo" export update:
#  - remove several rows that are no longer present in the refreshed export
#  - move the DANIELA (004329030) row so it sits right before the LUCYENE
#    (004376145) row, and update LUCYENE's balance
#
# All lookups are done by account number (column A) via Range.Find so the
# script does not depend on hard-coded row indices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Remove-AccountRow([string]$account) {
    $hit = $ws.Columns("A").Find($account, $null, $null, 1)
    if ($hit -ne $null) {
        $ws.Rows($hit.Row).Delete() | Out-Null
    }
}

# 1) Drop the duplicate-looking 004212581 / MARIA / 60000 row (004214460 row
#    with the same name stays untouched).
Remove-AccountRow "004212581"

# 2) Drop 004570632 / FABRICIO / 30490.28 entirely.
Remove-AccountRow "004570632"

# 3) Re-home 004329030 / DANIELA / 33665.17 so it comes right before the
#    004376145 / LUCYENE row, then drop its old location.
Remove-AccountRow "004329030"

$lucyene = $ws.Columns("A").Find("004376145", $null, $null, 1)
$ws.Rows($lucyene.Row).Insert() | Out-Null
$newRow = $lucyene.Row
$ws.Cells.Item($newRow, 1).Value = "'004329030"
$ws.Cells.Item($newRow, 2).Value = "DANIELA"
$ws.Cells.Item($newRow, 3).Value = 33665.17

# 4) Update LUCYENE's balance (row shifted down by one after the insert).
$lucyene2 = $ws.Columns("A").Find("004376145", $null, $null, 1)
$ws.Cells.Item($lucyene2.Row, 3).Value = 32000

# 5) Drop the remaining rows that were removed from the export.
Remove-AccountRow "005616259"
Remove-AccountRow "004813088"
Remove-AccountRow "005313179"
Remove-AccountRow "004455356"
Remove-AccountRow "004405234"
